$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eICR")
$ws.Range("A1").Value = "TEST"
